$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 2869.6667
$ws.Range("J100").Value = 2004
$ws.Range("L100").Value = 2004
$ws.Range("N100").Value = -3086
# Row 113
$ws.Range("H113").Value = 16111.111
$ws.Range("I113").Value = 17856.428
$ws.Range("K113").Value = 17856.428
$ws.Range("M113").Value = -14602.428
# Row 137
$ws.Range("H137").Value = 1642.2
$ws.Range("I137").Value = 1653.2858
$ws.Range("J137").Value = 1616.3334
$ws.Range("K137").Value = 4959.857400000001
$ws.Range("L137").Value = 4849.0002
$ws.Range("M137").Value = -2409.857400000001
$ws.Range("N137").Value = -9949.0002
# Row 138
$ws.Range("H138").Value = 3449.2307
$ws.Range("J138").Value = 4125
$ws.Range("L138").Value = 12375
$ws.Range("N138").Value = -22655

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 2485.889
$ws.Range("J88").Value = 4125
$ws.Range("L88").Value = 4125
$ws.Range("N88").Value = -4937
# Row 91
$ws.Range("H91").Value = 2485.889
$ws.Range("J91").Value = 4125
$ws.Range("L91").Value = 4125
$ws.Range("N91").Value = -6933
# Row 122
$ws.Range("H122").Value = 2539.1428
$ws.Range("I122").Value = 1199.5
$ws.Range("K122").Value = 3598.5
$ws.Range("M122").Value = -1148.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 900
$ws.Range("J12").Value = 900
$ws.Range("L12").Value = 900
$ws.Range("N12").Value = -1236
# Row 86
$ws.Range("H86").Value = 3059.0667
$ws.Range("I86").Value = 2643.6667
$ws.Range("K86").Value = 2643.6667
$ws.Range("M86").Value = -1520.6667
# Row 89
$ws.Range("H89").Value = 3059.0667
$ws.Range("I89").Value = 2643.6667
$ws.Range("K89").Value = 13218.3335
$ws.Range("M89").Value = -7602.333500000001
# Row 99
$ws.Range("H99").Value = 2163.6428
$ws.Range("I99").Value = 1899.2727
$ws.Range("J99").Value = 3133
$ws.Range("K99").Value = 1899.2727
$ws.Range("L99").Value = 3133
$ws.Range("M99").Value = -401.2727
$ws.Range("N99").Value = -6129

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 9941
$ws.Range("J28").Value = 9941
$ws.Range("L28").Value = 9941
$ws.Range("N28").Value = -10431
# Row 105
$ws.Range("H105").Value = 8663.615
$ws.Range("I105").Value = 10712.7
$ws.Range("J105").Value = 1833.3334
$ws.Range("K105").Value = 10712.7
$ws.Range("L105").Value = 1833.3334
$ws.Range("M105").Value = -8965.700000000001
$ws.Range("N105").Value = -5327.3334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 640
$ws.Range("I25").Value = 600
$ws.Range("J25").Value = 800
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2400
$ws.Range("M25").Value = -1631
$ws.Range("N25").Value = -2738
# Row 30
$ws.Range("H30").Value = 640
$ws.Range("I30").Value = 600
$ws.Range("J30").Value = 800
$ws.Range("K30").Value = 1800
$ws.Range("L30").Value = 2400
$ws.Range("M30").Value = -1698
$ws.Range("N30").Value = -2604
# Row 99
$ws.Range("H99").Value = 1380
$ws.Range("I99").Value = 1380
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4140
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1894
$ws.Range("N99").ClearContents()
# Row 129
$ws.Range("H129").Value = 2104.2
$ws.Range("I129").Value = 1089.4
$ws.Range("K129").Value = 3268.2
$ws.Range("M129").Value = 1731.8
# Row 130
$ws.Range("H130").Value = 17250
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 753.25
$ws.Range("I2").Value = 958.63635
$ws.Range("J2").Value = 502.22223
$ws.Range("K2").Value = 958.63635
$ws.Range("L2").Value = 502.22223
$ws.Range("M2").Value = -845.63635
$ws.Range("N2").Value = -728.2222300000001
# Row 80
$ws.Range("H80").Value = 2033.6
$ws.Range("I80").Value = 2061
$ws.Range("J80").Value = 1992.5
$ws.Range("K80").Value = 2061
$ws.Range("L80").Value = 1992.5
$ws.Range("M80").Value = -1063
$ws.Range("N80").Value = -3988.5
# Row 83
$ws.Range("H83").Value = 2033.6
$ws.Range("I83").Value = 2061
$ws.Range("J83").Value = 1992.5
$ws.Range("K83").Value = 10305
$ws.Range("L83").Value = 9962.5
$ws.Range("M83").Value = -5313
$ws.Range("N83").Value = -19946.5
# Row 113
$ws.Range("H113").Value = 1199
$ws.Range("I113").Value = 1199
$ws.Range("K113").Value = 1199
$ws.Range("M113").Value = 971
# Row 122
$ws.Range("H122").Value = 1774.2
$ws.Range("I122").Value = 1660.6
$ws.Range("J122").Value = 2001.4
$ws.Range("K122").Value = 4981.799999999999
$ws.Range("L122").Value = 6004.200000000001
$ws.Range("M122").Value = -2531.799999999999
$ws.Range("N122").Value = -10904.2
# Row 126
$ws.Range("H126").Value = 10380
$ws.Range("I126").Value = 8466.666999999999
$ws.Range("K126").Value = 25400.001
$ws.Range("M126").Value = -22930.001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 577.3125
$ws.Range("I55").Value = 139
$ws.Range("J55").Value = 639.9286
$ws.Range("K55").Value = 139
$ws.Range("L55").Value = 639.9286
$ws.Range("M55").Value = 34
$ws.Range("N55").Value = -985.9286
# Row 68
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251
# Row 71
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = -21256
# Row 93
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 933.3333
$ws.Range("K93").Value = 933.3333
$ws.Range("M93").Value = 314.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 168034.67
$ws.Range("J3").Value = 250052
$ws.Range("L3").Value = 250052
$ws.Range("N3").Value = -250280
# Row 81
$ws.Range("H81").Value = 5311.778
$ws.Range("I81").Value = 5404.4287
$ws.Range("J81").Value = 4987.5
$ws.Range("K81").Value = 10808.8574
$ws.Range("L81").Value = 9975
$ws.Range("M81").Value = -9747.857400000001
$ws.Range("N81").Value = -12097
# Row 84
$ws.Range("H84").Value = 5311.778
$ws.Range("I84").Value = 5404.4287
$ws.Range("J84").Value = 4987.5
$ws.Range("K84").Value = 54044.287
$ws.Range("L84").Value = 49875
$ws.Range("M84").Value = -48740.287
$ws.Range("N84").Value = -60483
# Row 96
$ws.Range("H96").Value = 1507.5834
$ws.Range("I96").Value = 1462.2106
$ws.Range("K96").Value = 1462.2106
$ws.Range("M96").Value = -89.21060000000011
# Row 100
$ws.Range("H100").Value = 3486442
$ws.Range("I100").Value = 4647323
$ws.Range("K100").Value = 9294646
$ws.Range("M100").Value = -9294105
# Row 126
$ws.Range("H126").Value = 3712.1667
$ws.Range("I126").Value = 3693.25
$ws.Range("K126").Value = 11079.75
$ws.Range("M126").Value = -8609.75
